$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''68.182.76'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  -0.45%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''3.885.78'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  -1.15%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("D5").Value = '''483.74'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  -0.13%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''145.29'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''  -2.35%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = '''0.621'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '''  -0.05%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = '''  -0.03%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = '''0.741'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''  +2.33%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = '''  +7.32%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = '''0.0000354'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '''  +0.80%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = '''43.31'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''  +1.72%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = '''  -0.51%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = '''4.508.68'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''  -1.13%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = '''3.889.88'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''  -2.96%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = '''14.20'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  -2.93%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = '''  -0.58%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = '''19.94'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''  +0.86%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = '''  +0.30%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = '''68.220.37'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''  -0.60%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = '''430.15'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''  -0.65%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = '''  +6.46%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = '''14.82'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''  +2.24%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = '''89.44'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''  +2.72%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = '''12.26'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '''  +18.29%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = '''3.69'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '''  +3.66%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = '''  -1.76%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = '''37.31'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '''  -2.35%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = '''  -3.83%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = '''711.50'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '''  -1.06%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = '''13.52'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '''  +1.94%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = '''  +0.33%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = '''  +2.91%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = '''0.0' + [char]0x2083 + '0880'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '''  -1.11%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = '''6.06'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '''  +9.75%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = '''  +3.59%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = '''41.06'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '''  -1.68%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("B38").Value = '''Kaspa'
$ws.Range("B38").Style = "Normal"
$ws.Range("C38").Value = '''https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("C38").Style = "Normal"
$ws.Range("D38").Value = '''0.145'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '''  -3.44%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("B39").Value = '''Fetch.AI'
$ws.Range("B39").Style = "Normal"
$ws.Range("C39").Value = '''https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("C39").Style = "Normal"
$ws.Range("D39").Value = '''3.02'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''  +5.68%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("B40").Value = '''Dai'
$ws.Range("B40").Style = "Normal"
$ws.Range("C40").Value = '''https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").Value = '''0.998'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''  -0.02%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("B41").Value = '''TheGraph'
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = '''https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = '''0.395'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''  +15.66%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = '''VeChain'
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = '''https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = '''0.0499'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '''  +6.36%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("E44").Value = '''  -1.95%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = '''  +1.44%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = '''3.35'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''  +3.90%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = '''  +0.10%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = '''3.37'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '''  -1.23%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = '''2.11'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '''  -2.51%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = '''144.47'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '''  -1.65%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = '''  -1.49%  '
$ws.Range("E51").Style = "Normal"
